$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 101
$ws.Range("B2").Value = "https://realpython.com/tutorials/web-dev/"
$ws.Range("C2").Value = 0.278
$ws.Range("D2").Value = 0.57
$ws.Range("E2").Value = 17
$ws.Range("F2").Value = 0.301
$ws.Range("G2").Value = 6.92
$ws.Range("H2").Value = 41
$ws.Range("I2").Value = 136
$ws.Range("J2").Value = 2.169
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 6.632

# Row 3
$ws.Range("A3").Value = 102
$ws.Range("B3").Value = "https://python.land/introduction-to-python"
$ws.Range("C3").Value = 0.19
$ws.Range("D3").Value = 0.47
$ws.Range("E3").Value = 8.5
$ws.Range("F3").Value = 0.206
$ws.Range("G3").Value = 3.482
$ws.Range("H3").Value = 14
$ws.Range("I3").Value = 68
$ws.Range("J3").Value = 2.015
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 6.206

# Row 4 (new row)
$ws.Range("A4").Value = 103
$ws.Range("B4").Value = "https://en.wikipedia.org/wiki/Natural_language_processing"
$ws.Range("C4").Value = 0.03
$ws.Range("D4").Value = 0.429
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 0.419
$ws.Range("G4").Value = 6.168
$ws.Range("H4").Value = 673
$ws.Range("I4").Value = 1605
$ws.Range("J4").Value = 2.414
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 7.128
